$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.294.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.41%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.656.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.45%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.67%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.57%  '

$ws.Range('E9').Value = '  +9.17%  '

$ws.Range('E10').Value = '  +4.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.61%  '

$ws.Range('E12').Value = '  +1.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000190'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +16.67%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.136.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.59%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.129.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.40%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.656.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.63%  '

$ws.Range('E22').Value = '  -0.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.86%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.66%  '

$ws.Range('E26').Value = '  +18.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.61%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.47%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.165'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.56%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.51%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '549.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.31%  '

$ws.Range('E32').Value = '  +0.08%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.46%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.76%  '

$ws.Range('E36').Value = '  +3.99%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.74%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.22%  '

$ws.Range('E40').Value = '  +0.07%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '167.96'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.66%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0622'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.47%  '

$ws.Range('E46').Value = '  +8.67%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.69%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.661'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0265'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0986'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.18%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.87'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.61%  '
